# 15.1.2.1.xlsx update:
#  - extend the data table with two more years (2022, 2023)
#  - matching value 6.53 for the new "Protected areas" row
#  - adjust a few row heights to the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (year headers): add N4=2022, O4=2023 --------------------------
# Copy M4 (style s="4") into N4/O4 so the new cells inherit the same
# number format / font / border, then overwrite with the real year values.
$ws.Cells.Item(4, 13).Copy($ws.Cells.Item(4, 14))
$ws.Cells.Item(4, 13).Copy($ws.Cells.Item(4, 15))
$ws.Cells.Item(4, 14).Value = 2022
$ws.Cells.Item(4, 15).Value = 2023

# --- Row 5 (values): add N5=6.53, O5=6.53 ---------------------------------
# Copy M5 (style s="5", already holding 6.53) into N5/O5 so both format and
# the 6.53 figure carry over exactly.
$ws.Cells.Item(5, 13).Copy($ws.Cells.Item(5, 14))
$ws.Cells.Item(5, 13).Copy($ws.Cells.Item(5, 15))

# --- Row height tweaks -----------------------------------------------------
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 17.25
